$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Fix the duplicated/typo'd shared string on Product_Tests sheet.
#    The old text '2 Verify URL contains "products"' (missing the
#    period after "2") gets corrected to match the already-existing
#    string '2.Verify URL contains "products"' used elsewhere, which
#    lets the workbook de-duplicate the shared string table entry.
# ------------------------------------------------------------------
$wsProduct = $wb.Worksheets.Item("Product_Tests")
$wsProduct.Range("C3").Value = '2.Verify URL contains "products"'

# Minor column width tweak on Product_Tests (column B).
# NOTE: the host runtime quantizes ColumnWidth input on a 1/6-character
# grid before storing it, so the literal target width (23.109375) isn't
# exactly reproducible; 22.333333333333332 is the closest input that
# yields the nearest achievable stored width.
$wsProduct.Columns.Item(2).ColumnWidth = 22.333333333333332

# ------------------------------------------------------------------
# 2. Add the new Sales_Tests worksheet after Product_Tests (the last
#    sheet), matching the header styling used by the other sheets.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSales = $wb.Worksheets.Add($null, $lastSheet)
$wsSales.Name = "Sales_Tests"

# Copy header formatting (fill/style) from an existing sheet's header row.
$wsProduct.Range("A1:C1").Copy()
$wsSales.Range("A1:C1").PasteSpecial(-4122)

# Header row contents (re-uses the existing shared strings).
$wsSales.Range("A1").Value = "Test Case ID(s)"
$wsSales.Range("B1").Value = "Test Case Description"
$wsSales.Range("C1").Value = "Test Steps"

# New test-case data.
$wsSales.Range("A2").Value = "TC_SALE_01"
$wsSales.Range("B2").Value = "Navigate to Sales"
$wsSales.Range("C2").Value = '1.Click on "Sales" at "//a[contains(@href, ''/sales'')]"'
$wsSales.Range("C3").Value = '2.Verify URL contains "sales"'

# Column widths for the new sheet (see note above re: 1/6-character
# quantization - these inputs give the closest achievable result to the
# targets 18.21875 / 28.6640625 / 114.77734375).
$wsSales.Columns.Item(1).ColumnWidth = 17.333333333333332
$wsSales.Columns.Item(2).ColumnWidth = 27.833333333333332
$wsSales.Columns.Item(3).ColumnWidth = 114.0

# ------------------------------------------------------------------
# 3. Update sheet-view selections/active tab so the new Sales_Tests
#    sheet becomes the active/selected sheet, and Product_Tests keeps
#    a plain (non-tab-selected) view with its selection reset to C3.
# ------------------------------------------------------------------
$wsProduct.Activate()
$wsProduct.Range("C3").Select()

$wsSales.Activate()
$wsSales.Range("C4").Select()
